# Updates the cryptos list: refresh Price (column D) and Volume(1h) (column E)
# for the rows whose quotes moved since the last scrape.
#
# Column D values are stored as plain text (e.g. "28.421.08", "0.9997") even
# though most of them look numeric. Assigning them straight to .Value would
# let Excel auto-coerce them into real numbers (losing trailing zeros /
# exact formatting, e.g. "1.000" -> 1). Prefixing with a single quote forces
# Excel to keep the literal text, and resetting .Style back to "Normal"
# afterwards clears the "quote prefix" flag so no stray cell style lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.421.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.24%  '
$ws.Range("D3").Value = "'1.803.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'314.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = "'0.5516"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.69%  '
$ws.Range("D8").Value = "'0.3860"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.87%  '
$ws.Range("D9").Value = "'0.07608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.28%  '
$ws.Range("D10").Value = "'42.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").Value = "'1.129"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.51%  '
$ws.Range("D12").Value = "'0.9998"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("E13").Value = '  +3.85%  '
$ws.Range("D14").Value = "'6.189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("D15").Value = "'7.436"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.07%  '
$ws.Range("D16").Value = "'1.806.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'17.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.78%  '
$ws.Range("D22").Value = "'5.976"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.58%  '
$ws.Range("D23").Value = "'28.437.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.00%  '
$ws.Range("D24").Value = "'11.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("D25").Value = "'2.142"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.10%  '
$ws.Range("D26").Value = "'159.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").Value = "'20.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.08%  '
$ws.Range("D28").Value = "'2.416"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.80%  '
$ws.Range("D29").Value = "'2.014.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("D30").Value = "'123.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").Value = "'1.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.85%  '
$ws.Range("D32").Value = "'0.1024"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.98%  '
$ws.Range("D33").Value = "'5.773"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.78%  '
$ws.Range("D34").Value = "'3.688"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").Value = "'0.2316"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +14.67%  '
$ws.Range("D36").Value = "'0.06431"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.90%  '
$ws.Range("D37").Value = "'0.02330"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.59%  '
$ws.Range("E38").Value = '  +7.05%  '
$ws.Range("D39").Value = "'8.810"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.24%  '
$ws.Range("D40").Value = "'11.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.00%  '
$ws.Range("D41").Value = "'0.6421"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.87%  '
$ws.Range("D42").Value = "'1.162"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.17%  '
$ws.Range("D43").Value = "'0.9992"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = '  -3.35%  '
$ws.Range("D45").Value = "'13.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.43%  '
$ws.Range("D46").Value = "'0.5990"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.28%  '
$ws.Range("D47").Value = "'3.683"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("D48").Value = "'127.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.22%  '
$ws.Range("D49").Value = "'1.986"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = '  +3.60%  '
$ws.Range("D51").Value = "'0.06899"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.77%  '
